$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: remove the "Meta description: ..." paragraph that
# currently follows the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# Change 2: insert a new paragraph right before the very last
# paragraph (the one with the italic "feature image" prompt),
# containing an empty run followed by a bold run with the title text.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$insertPos = $secondLast.Range.End
$insertRange = $d.Range($insertPos, $insertPos)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Clover Lady Free Today! Review &amp; Ratings</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRange.InsertXML($xml)

# InsertXML leaves behind a stray empty paragraph used to force the
# paragraph break; remove it.
$strayPara = $d.Paragraphs.Item($count + 1)
$strayPara.Range.Delete()

# ------------------------------------------------------------------
# Change 3: replace the text of the final paragraph's italic run
# (the old "feature image" art-direction prompt) with the meta
# description text, keeping its italic formatting intact.
# ------------------------------------------------------------------
$finalCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($finalCount)
$lastPara.Range.Find.Execute(
    "For the feature image of Clover Lady, let's have a cartoon-style Maya warrior with glasses. The image should feature the Maya warrior happily playing the game on a mobile device or computer, with the magical forest and mushroom-shaped game grid in the background. The warrior should be holding a clover symbol, with the Metalwolf and girl bonus symbols also visible. The overall style should be colorful and playful, capturing the fairy tale theme of the game.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Explore the enchanting forest with Clover Lady. Read the review, play for free, and discover bonus features, graphics, and design. Compatible on all devices.",
    2)
